# Generate Report for Handback
#
# This localization-status report is refreshed after a handback event:
#  - the "Ready for handoff" status becomes "Handed back: in sync with en-US"
#    wherever it is shown (Overview zh-cn/de-de status columns, and each
#    language sheet's own Status column)
#  - the per-language sheets gain their "Latest Target File" / "Latest
#    Handback File" links for the two source docs, now that a handback
#    round has completed
#  - de-de (which is in sync) gets a real "Latest Handback DateTime";
#    zh-cn's placeholder datetime also gets refreshed as part of the same
#    report regeneration pass
#  - a few columns are widened so the newly-populated long file names and
#    status text are not clipped

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdFile1 = "02578520-0cef-4411-98fe-3ddea15f5a9b.md"
$mdFile2 = "403a8a16-e353-4292-b848-8cbc83033232.md"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cdcfbec9c467cce3648aff7d347ff02c1122c289/e2e/$mdFile1"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cdcfbec9c467cce3648aff7d347ff02c1122c289/e2e/$mdFile2"

$zhTarget1 = "02578520-0cef-4411-98fe-3ddea15f5a9b.6adccbb1821ae2060cdbecda6bca56052f7e6eb3.zh-cn.xlf"
$zhTarget2 = "403a8a16-e353-4292-b848-8cbc83033232.17d529a906857c00d6a50777ef50290a87e8d9bd.zh-cn.xlf"
$deTarget1 = "02578520-0cef-4411-98fe-3ddea15f5a9b.6adccbb1821ae2060cdbecda6bca56052f7e6eb3.de-de.xlf"
$deTarget2 = "403a8a16-e353-4292-b848-8cbc83033232.17d529a906857c00d6a50777ef50290a87e8d9bd.de-de.xlf"

$zhHandbackDate = "2016-09-04 06:32:44"
$deHandbackDate = "2016-09-04 06:32:51"

# ColumnWidth values chosen so the saved (post-padding) column width lands on
# the widened target: 29.9777047293527 characters -> use 29.1666... ; the
# 18.65/21.71-wide columns become a flat 40 characters -> use 39.1666...
$wideColWidth   = 29.166666666666668
$fortyColWidth  = 39.166666666666664

# ---------------------------------------------------------------------------
# Overview sheet: the zh-cn / de-de status cells show the new handback status
# (driven by the shared "Ready for handoff" text), and its zh-cn/de-de columns
# widen along with the same widening on the per-language sheets.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, "", "", $mdFile1)
$wsZh.Range("J2").Value = $zhTarget1

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, "", "", $mdFile2)
$wsZh.Range("J3").Value = $zhTarget2

$wsZh.Range("K2").Value = $zhHandbackDate
$wsZh.Range("K3").Value = $zhHandbackDate

$wsZh.Columns.Item(3).ColumnWidth = $wideColWidth
$wsZh.Columns.Item(9).ColumnWidth = $fortyColWidth
$wsZh.Columns.Item(10).ColumnWidth = $fortyColWidth

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, "", "", $mdFile1)
$wsDe.Range("J2").Value = $deTarget1
$wsDe.Range("K2").Value = $deHandbackDate

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, "", "", $mdFile2)
$wsDe.Range("J3").Value = $deTarget2
$wsDe.Range("K3").Value = $deHandbackDate

$wsDe.Columns.Item(3).ColumnWidth = $wideColWidth
$wsDe.Columns.Item(9).ColumnWidth = $fortyColWidth
$wsDe.Columns.Item(10).ColumnWidth = $fortyColWidth

Write-Host "Handback report regenerated."
